$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$refStyle = $ws.Range("A517").Style

$data = @(
    @('000517', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_1_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_1_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 07-Jun-2023 12:59:42'),
    @('000518', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_2_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 07-Jun-2023 12:59:43'),
    @('000519', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 07-Jun-2023 12:59:43'),
    @('000520', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:03:58'),
    @('000521', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:04:14'),
    @('000522', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:04:29'),
    @('000523', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:04:33'),
    @('000524', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:04:39'),
    @('000525', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:06:03'),
    @('000526', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:06:07'),
    @('000527', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:15:50'),
    @('000528', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:16:08'),
    @('000529', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:16:18'),
    @('000530', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:17:17'),
    @('000531', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. Script used: ResampledVisualisationCycles_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_Exp_B1_001_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:17:37'),
    @('000532', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Waist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_C1_001_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:19:07'),
    @('000533', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Knees. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_B1_001_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:20:25'),
    @('000534', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Knees. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:26:39'),
    @('000535', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Knees. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:26:54'),
    @('000536', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_D1_002_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:27:04'),
    @('000537', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Feet. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_E1_001_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:27:49'),
    @('000538', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Wrist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_E1_001_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:28:11'),
    @('000539', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Wrist. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_E1_001_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:28:54'),
    @('000540', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Upper_Arm. Script used: ResampledVisualisationCycles_XY.  Dataset used: Session01_Exp_F1_001_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_IndividualCyclesTimeThreshold. Date Generated: 07-Jun-2023 13:29:05'),
    @('000541', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 13:30:07'),
    @('000542', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 13:58:56'),
    @('000543', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 13:59:38'),
    @('000544', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:03:29'),
    @('000545', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:04:08'),
    @('000546', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:04:57'),
    @('000547', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:07:46'),
    @('000548', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:08:27'),
    @('000549', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:08:49'),
    @('000550', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:10:47'),
    @('000551', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:12:27'),
    @('000552', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:19:36'),
    @('000553', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:19:56'),
    @('000554', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:20:22'),
    @('000555', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:20:36'),
    @('000556', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:20:48'),
    @('000557', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:21:35'),
    @('000558', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:21:54'),
    @('000559', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:29:09'),
    @('000560', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp E1 001 SavedCycle 2. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:30:03'),
    @('000561', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp D1 001 SavedCycle 2. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:31:52'),
    @('000562', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp D1 002 SavedCycle 2. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:32:07'),
    @('000563', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp F1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:32:55'),
    @('000564', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp C1 001 SavedCycle 3. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:33:37'),
    @('000565', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp C1 001 SavedCycle 4. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:33:51'),
    @('000566', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp B1 001 SavedCycle 2. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:34:23'),
    @('000567', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp B1 001 SavedCycle 2. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:34:53'),
    @('000568', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp A1 003 SavedCycle 2. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:36:23'),
    @('000569', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp A1 003 SavedCycle 4. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:36:30'),
    @('000570', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp A1 005 SavedCycle 4. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:36:37'),
    @('000571', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath. Arrays shifted to start at highest point.. Script used: RigidRegistrationAndScalingWithDeviation.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01 Exp B1 001 SavedCycle 2. File Location: Visualisations/RegistrationAnalysis. Date Generated: 07-Jun-2023 14:37:35')
)

$startRow = 518
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $guid = $data[$idx][0]
    $info = $data[$idx][1]
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $aCell.NumberFormat = "@"
    $aCell.Value = $guid
    $aCell.Style = $refStyle
    $bCell.Value = $info
    $bCell.Style = $refStyle
}
